$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 66, shifting existing rows 66-176 down to 67-177.
$ws.Rows(66).Insert()

# Populate the newly inserted row 66 with the new record's data.
$ws.Range("A66").Value = 5
$ws.Range("B66").Value = "Macroferia Regional de Talca"
$ws.Range("C66").Value = "Maule"
$ws.Range("D66").Value = 44469
$ws.Range("D66").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E66").Value = 7
$ws.Range("F66").Value = 100112009
$ws.Range("G66").Value = "Acelga"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 1800
$ws.Range("L66").Value = 1800
$ws.Range("M66").Value = 1800
$ws.Range("N66").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O66").Value = "Región del Maule"
$ws.Range("P66").Value = 450
$ws.Range("Q66").Value = 4
$ws.Range("R66").Value = "Hortaliza"
